$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Typo fix: "Burma CSS" -> "Bulma CSS" in the Technical Skills paragraph.
#    The edit also re-splits the surrounding text into three runs (matching
#    how Word splits a run when you select/retype a word) and moves the
#    "_GoBack" bookmark to sit right after the corrected word "Bulma".
# ---------------------------------------------------------------------------
$rTypo = $d.Content
$null = $rTypo.Find.Execute("Burma CSS", $true, $false, $false, $false, $false, $true, 1, $false, "Bulma CSS", 2)

$rBulma = $d.Content
$null = $rBulma.Find.Execute("Bulma", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$bulmaStart = $rBulma.Start
$bulmaEnd = $rBulma.End

# Force a run split right before "Bulma" (temporary bookmark), so "Bootstrap, "
# and "Bulma" end up as distinct runs instead of being re-merged.
$splitPoint = $d.Range($bulmaStart, $bulmaStart)
$d.Bookmarks.Add("TempSplit", $splitPoint)

# Move the "_GoBack" bookmark from its old spot (end of "Key Accomplishments:")
# to sit immediately after "Bulma".
$oldGoBack = $d.Bookmarks("_GoBack")
$oldGoBack.Delete()

$afterBulma = $d.Range($bulmaEnd, $bulmaEnd)
$d.Bookmarks.Add("_GoBack", $afterBulma)

$tempBm = $d.Bookmarks("TempSplit")
$tempBm.Delete()

# ---------------------------------------------------------------------------
# 2) Merge the " | " run and "Deployed:" run into a single " | Deployed:" run
#    (first occurrence only - the GitHub / Deployed line).
# ---------------------------------------------------------------------------
$rAnchor = $d.Content
$null = $rAnchor.Find.Execute("3zXw7W1", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$afterLink = $rAnchor.End

$rDepSearch = $d.Range($afterLink, $d.Content.End)
$null = $rDepSearch.Find.Execute(" | Deployed:", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$depStart = $rDepSearch.Start
$depEnd = $rDepSearch.End

$pipeLen = 3  # " | "
$rPipe = $d.Range($depStart, $depStart + $pipeLen)
$rDeployedWord = $d.Range($depStart + $pipeLen, $depEnd)
$deployedText = $rDeployedWord.Text
$rDeployedWord.Text = ""

$rPipeAgain = $d.Range($depStart, $depStart + $pipeLen)
$rPipeAgain.InsertAfter($deployedText)

# ---------------------------------------------------------------------------
# 3) Merge "Created ", "the entire" and " front end..." runs into a single run.
# ---------------------------------------------------------------------------
$rCreatedFull = $d.Content
$null = $rCreatedFull.Find.Execute("Created the entire front end HTML/CSS and JavaScript for DOM manipulation and connection and consumption of APIs.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$cStart = $rCreatedFull.Start
$cEnd = $rCreatedFull.End

$prefixLen = 18  # "Created the entire"
$rPrefix = $d.Range($cStart, $cStart + $prefixLen)
$prefixText = $rPrefix.Text
$rPrefix.Text = ""

$rRemainder = $d.Range($cStart, $cStart + ($cEnd - $cStart - $prefixLen))
$rRemainder.InsertBefore($prefixText)

# Nudge the merged run so the engine recomputes xml:space (drops the stray
# "preserve" left over from the original runs, matching a clean merged run).
$rMerged = $d.Range($cStart, $cEnd)
$mergedText = $rMerged.Text
$rMerged.Text = $mergedText + "X"
$rMergedPlusX = $d.Range($cStart, $cEnd + 1)
$rMergedPlusX.Text = $mergedText
